# Applies the "dEntrevista" -> "Entrevista 1" title fix and merges the
# split "-Quais ti" / "pos de " runs in the last paragraph into a single
# run reading "-Quais tipos de ".

$d = $word.ActiveDocument

# 1) Fix the title paragraph: "dEntrevista" -> "Entrevista 1"
$d.Content.Find.Execute("dEntrevista", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Entrevista 1", 2)

# 2) Merge "-Quais ti" + "pos de " into "-Quais tipos de "
$d.Content.Find.Execute("-Quais tipos de ", $false, $false, $false, $false, $false,
                         $true, 1, $false, "-Quais tipos de ", 2)

$d.Save()
